$d = $word.ActiveDocument

# The table's first two rows (header labels row + the "n" row) get their
# cells' vertical alignment set to bottom (w:vAlign w:val="bottom").
$tbl = $d.Tables.Item(1)

for ($r = 1; $r -le 2; $r++) {
    $row = $tbl.Rows.Item($r)
    foreach ($cell in $row.Cells) {
        $cell.VerticalAlignment = 3
    }
}
